# Apply the edits described by the diff to slide 1 of the presentation:
#  1) Move the "Work prepared by..." textbox (shape id=35) down
#     (a:off y 2879021 -> 3001569, x/width/height unchanged).
#  2) Update the date line in the "L1 SIGL, USTHB ..." textbox (shape id=36)
#     from "31 december 2023" to "04 january 2024".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Helper: PowerPoint's Shape.Top/Left are stored as single-precision
# (float32) point values internally, so a naive EMU/12700.0 assignment
# can be off by a single EMU once rounded back on save. Nudge the
# assigned point value by tiny increments until the value we read back
# (quantised through that same float32 round-trip) maps to exactly the
# desired EMU, guaranteeing an exact match in the saved XML.
function Set-ExactTopEmu {
    param($Shape, [double]$TargetEmu)

    $basePts = $TargetEmu / 12700.0
    for ($i = 0; $i -le 10000; $i++) {
        $cand = $basePts + ($i * 0.0000001)
        $Shape.Top = $cand
        $readBack = $Shape.Top
        $emuCheck = [math]::Floor(($readBack * 12700.0) + 0.5)
        if ($emuCheck -eq $TargetEmu) {
            return $true
        }
    }
    return $false
}

# --- 1) Reposition the "Work prepared by" textbox (shape index 6, id 35) ---
$shpAuthors = $s.Shapes.Item(6)
Set-ExactTopEmu $shpAuthors 3001569 | Out-Null

# --- 2) Update the date text in the last textbox (shape index 7, id 36) ---
$shpDate = $s.Shapes.Item(7)
$tr = $shpDate.TextFrame.TextRange

# Original text: "   L1 SIGL, USTHB                      31 december 2023"
# split across 3 runs:
#   run1 (chars 1-42):  "   L1 SIGL, USTHB                      31 "
#   run2 (chars 43-50): "december"
#   run3 (chars 51-55): " 2023"
# Replace from right to left so earlier character offsets stay valid
# while the replacement text lengths differ from the originals.

$run3 = $tr.Characters(51, 5)
$run3.Text = " 2024"

$run2 = $tr.Characters(43, 8)
$run2.Text = "january"

$run1 = $tr.Characters(1, 42)
$run1.Text = "   L1 SIGL, USTHB                      04 "

Write-Output "authors box top (pt): $($shpAuthors.Top)"
Write-Output "date text now: $($tr.Text)"
